$d = $word.ActiveDocument

$firstPara = $d.Paragraphs.First.Range
$firstPara.InsertBefore("Release Notes`r")

$headingRange = $d.Paragraphs.First.Range
$headingRange.Style = "Heading 1"
